$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly cryptocurrency price/volume refresh (coinranking.com data)

# Step 1: mark numeric-looking Price cells as Text so values like "1.006"
# keep their exact literal formatting (matches source data which stores
# prices as plain text, not numbers).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Step 2: write the refreshed values

# Row 2
$ws.Range('D2').Value = '28.385.42'
$ws.Range('E2').Value = '  -2.85%  '

# Row 3
$ws.Range('D3').Value = '1.954.18'
$ws.Range('E3').Value = '  -0.50%  '

# Row 4
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  -0.90%  '

# Row 5
$ws.Range('D5').Value = '320.01'
$ws.Range('E5').Value = '  -2.55%  '

# Row 6
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').Value = '  -1.05%  '

# Row 7
$ws.Range('D7').Value = '0.4774'
$ws.Range('E7').Value = '  -4.25%  '

# Row 8
$ws.Range('D8').Value = '0.4026'

# Row 9
$ws.Range('D9').Value = '53.73'
$ws.Range('E9').Value = '  +0.27%  '

# Row 10
$ws.Range('D10').Value = '0.08456'
$ws.Range('E10').Value = '  -5.02%  '

# Row 11
$ws.Range('D11').Value = '1.053'
$ws.Range('E11').Value = '  -4.31%  '

# Row 12
$ws.Range('D12').Value = '22.33'
$ws.Range('E12').Value = '  -2.98%  '

# Row 13
$ws.Range('D13').Value = '1.950.02'
$ws.Range('E13').Value = '  -1.01%  '

# Row 14
$ws.Range('D14').Value = '7.550'
$ws.Range('E14').Value = '  -4.08%  '

# Row 15
$ws.Range('D15').Value = '6.144'
$ws.Range('E15').Value = '  -4.10%  '

# Row 16
$ws.Range('D16').Value = '1.006'
$ws.Range('E16').Value = '  -1.03%  '

# Row 17
$ws.Range('D17').Value = '90.74'
$ws.Range('E17').Value = '  -1.03%  '

# Row 18
$ws.Range('D18').Value = '0.00001065'
$ws.Range('E18').Value = '  -3.40%  '

# Row 19
$ws.Range('D19').Value = '0.06583'
$ws.Range('E19').Value = '  -1.92%  '

# Row 20
$ws.Range('D20').Value = '18.47'
$ws.Range('E20').Value = '  -4.20%  '

# Row 21
$ws.Range('E21').Value = '  -1.13%  '

# Row 22
$ws.Range('D22').Value = '5.816'
$ws.Range('E22').Value = '  -1.87%  '

# Row 23
$ws.Range('D23').Value = '28.388.44'
$ws.Range('E23').Value = '  -2.95%  '

# Row 24
$ws.Range('E24').Value = '  -3.77%  '

# Row 25
$ws.Range('E25').Value = '  -1.17%  '

# Row 26
$ws.Range('D26').Value = '2.186.22'
$ws.Range('E26').Value = '  -1.79%  '

# Row 27
$ws.Range('D27').Value = '154.48'
$ws.Range('E27').Value = '  -0.62%  '

# Row 28
$ws.Range('D28').Value = '20.20'
$ws.Range('E28').Value = '  -2.09%  '

# Row 29
$ws.Range('D29').Value = '5.913'
$ws.Range('E29').Value = '  -4.82%  '

# Row 30
$ws.Range('D30').Value = '2.154'
$ws.Range('E30').Value = '  -6.10%  '

# Row 31
$ws.Range('D31').Value = '123.87'
$ws.Range('E31').Value = '  -2.34%  '

# Row 32
$ws.Range('D32').Value = '0.9775'
$ws.Range('E32').Value = '  -6.99%  '

# Row 33
$ws.Range('D33').Value = '0.09622'
$ws.Range('E33').Value = '  -2.56%  '

# Row 34
$ws.Range('D34').Value = '1.451'
$ws.Range('E34').Value = '  -4.23%  '

# Row 35
$ws.Range('D35').Value = '5.601'
$ws.Range('E35').Value = '  -3.29%  '

# Row 36
$ws.Range('E36').Value = '  -2.28%  '

# Row 37
$ws.Range('D37').Value = '8.960'
$ws.Range('E37').Value = '  -2.08%  '

# Row 38
$ws.Range('D38').Value = '0.02329'
$ws.Range('E38').Value = '  -4.05%  '

# Row 39
$ws.Range('D39').Value = '0.06220'
$ws.Range('E39').Value = '  -1.75%  '

# Row 40
$ws.Range('D40').Value = '1.250'
$ws.Range('E40').Value = '  -2.94%  '

# Row 41
$ws.Range('D41').Value = '0.6205'
$ws.Range('E41').Value = '  -4.23%  '

# Row 42
$ws.Range('D42').Value = '11.13'
$ws.Range('E42').Value = '  -3.49%  '

# Row 43
$ws.Range('D43').Value = '1.004'
$ws.Range('E43').Value = '  -1.09%  '

# Row 44
$ws.Range('D44').Value = '0.1917'
$ws.Range('E44').Value = '  -5.10%  '

# Row 45
$ws.Range('D45').Value = '1.347'
$ws.Range('E45').Value = '  +5.00%  '

# Row 46
$ws.Range('D46').Value = '0.5960'
$ws.Range('E46').Value = '  -4.77%  '

# Row 47
$ws.Range('E47').Value = '  -3.39%  '

# Row 48
$ws.Range('D48').Value = '2.060'
$ws.Range('E48').Value = '  -5.75%  '

# Row 49
$ws.Range('D49').Value = '3.382'
$ws.Range('E49').Value = '  -2.96%  '

# Row 50
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.06802'
$ws.Range('E50').Value = '  -0.79%  '

# Row 51
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.00000000314'
$ws.Range('E51').Value = '  -7.63%  '
